$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New attendance log row (row 2) appended under the existing header row.
# Date / ClockIn look like dates/times to Excel's smart-typing, so we
# briefly mark the cell as Text before writing the literal string, then
# drop back to the Normal style so no stray number-format sticks around.

$a2 = $ws.Range("A2")
$a2.NumberFormat = "@"
$a2.Value = "10/05/2025"
$a2.Style = "Normal"

$ws.Range("B2").Value = 2602069620

$c2 = $ws.Range("C2")
$c2.NumberFormat = "@"
$c2.Value = "14:05:33"
$c2.Style = "Normal"

# ClockOut / Log are blank for this row, but still present as empty cells.
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
